$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 353.55554
$ws.Range("J2").Value = 330.33334
$ws.Range("L2").Value = 330.33334
$ws.Range("N2").Value = -556.33334
$ws.Range("H11").Value = 6.923077
$ws.Range("I11").Value = 6.923077
$ws.Range("K11").Value = 6.923077
$ws.Range("M11").Value = 133.076923
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H40").Value = 1252.2858
$ws.Range("I40").Value = 759.1667
$ws.Range("J40").Value = 1622.125
$ws.Range("K40").Value = 759.1667
$ws.Range("L40").Value = 1622.125
$ws.Range("M40").Value = -584.1667
$ws.Range("N40").Value = -1972.125
$ws.Range("H45").Value = 1100
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 1100
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 3300
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -3684
$ws.Range("H51").Value = 3227
$ws.Range("J51").Value = 2395.1
$ws.Range("L51").Value = 2395.1
$ws.Range("N51").Value = -3363.1
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H132").Value = 3218.25
$ws.Range("I132").Value = 3869.7083
$ws.Range("K132").Value = 11609.1249
$ws.Range("M132").Value = -9079.124899999999
$ws.Range("H137").Value = 1261.919
$ws.Range("I137").Value = 1255.4231
$ws.Range("J137").Value = 1277.2727
$ws.Range("K137").Value = 3766.2693
$ws.Range("L137").Value = 3831.8181
$ws.Range("M137").Value = -1216.2693
$ws.Range("N137").Value = -8931.8181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4130.75
$ws.Range("I61").Value = 4163.4375
$ws.Range("K61").Value = 4163.4375
$ws.Range("M61").Value = -3951.4375
$ws.Range("H74").Value = 29413336
$ws.Range("I74").Value = 47619700
$ws.Range("K74").Value = 47619700
$ws.Range("M74").Value = -47618826
$ws.Range("H77").Value = 29413336
$ws.Range("I77").Value = 47619700
$ws.Range("K77").Value = 238098500
$ws.Range("M77").Value = -238094132
$ws.Range("H136").Value = 4130.75
$ws.Range("I136").Value = 4163.4375
$ws.Range("K136").Value = 12490.3125
$ws.Range("M136").Value = -9940.3125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 690.5
$ws.Range("I22").Value = 640.7273
$ws.Range("J22").Value = 800
$ws.Range("K22").Value = 640.7273
$ws.Range("L22").Value = 800
$ws.Range("M22").Value = -467.7273
$ws.Range("N22").Value = -1146

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 4000
$ws.Range("J4").Value = 4000
$ws.Range("L4").Value = 4000
$ws.Range("N4").Value = -4224
$ws.Range("H22").Value = 183
$ws.Range("I22").Value = 183
$ws.Range("K22").Value = 183
$ws.Range("M22").Value = 167
$ws.Range("H23").Value = 3800
$ws.Range("J23").Value = 3800
$ws.Range("L23").Value = 3800
$ws.Range("N23").Value = -4280
$ws.Range("H27").Value = 3800
$ws.Range("J27").Value = 3800
$ws.Range("L27").Value = 3800
$ws.Range("N27").Value = -4184
$ws.Range("H132").Value = 6112.4
$ws.Range("I132").Value = 5622.4
$ws.Range("J132").Value = 6602.4
$ws.Range("K132").Value = 16867.2
$ws.Range("L132").Value = 19807.2
$ws.Range("M132").Value = -14337.2
$ws.Range("N132").Value = -24867.2
$ws.Range("H134").Value = 1449.9286
$ws.Range("I134").Value = 1199.8889
$ws.Range("K134").Value = 3599.6667
$ws.Range("M134").Value = -1064.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H119").Value = 6770
$ws.Range("I119").Value = 5155
$ws.Range("J119").Value = 10000
$ws.Range("K119").Value = 15465
$ws.Range("L119").Value = 30000
$ws.Range("M119").Value = -10627
$ws.Range("N119").Value = -39676
$ws.Range("H131").Value = 696.5599999999999
$ws.Range("J131").Value = 696.5599999999999
$ws.Range("L131").Value = 2089.68
$ws.Range("N131").Value = -12169.68

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 8363.637000000001
$ws.Range("J5").Value = 8900
$ws.Range("L5").Value = 8900
$ws.Range("N5").Value = -9124
$ws.Range("H51").Value = 30326
$ws.Range("J51").Value = 30326
$ws.Range("L51").Value = 30326
$ws.Range("N51").Value = -31344
$ws.Range("H102").Value = 1494.875
$ws.Range("I102").Value = 1203.0344
$ws.Range("K102").Value = 1203.0344
$ws.Range("M102").Value = 418.9656
$ws.Range("H126").Value = 2852.3713
$ws.Range("I126").Value = 2102.4092
$ws.Range("J126").Value = 4121.5386
$ws.Range("K126").Value = 6307.2276
$ws.Range("L126").Value = 12364.6158
$ws.Range("M126").Value = -3837.2276
$ws.Range("N126").Value = -17304.6158
$ws.Range("H136").Value = 8249.5
$ws.Range("J136").Value = 8249.5
$ws.Range("L136").Value = 24748.5
$ws.Range("N136").Value = -29848.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 241158.4
$ws.Range("I2").Value = 420000.4
$ws.Range("J2").Value = 62316.4
$ws.Range("K2").Value = 420000.4
$ws.Range("L2").Value = 62316.4
$ws.Range("M2").Value = -419888.4
$ws.Range("N2").Value = -62540.4
$ws.Range("H22").Value = 4790.1
$ws.Range("I22").Value = 5650.125
$ws.Range("J22").Value = 1350
$ws.Range("K22").Value = 5650.125
$ws.Range("L22").Value = 1350
$ws.Range("M22").Value = -5355.125
$ws.Range("N22").Value = -1940
$ws.Range("H25").Value = 4000
$ws.Range("J25").Value = 4000
$ws.Range("L25").Value = 4000
$ws.Range("N25").Value = -4460
$ws.Range("H27").Value = 4790.1
$ws.Range("I27").Value = 5650.125
$ws.Range("J27").Value = 1350
$ws.Range("K27").Value = 5650.125
$ws.Range("L27").Value = 1350
$ws.Range("M27").Value = -5543.125
$ws.Range("N27").Value = -1564
$ws.Range("H132").Value = 504065.25
$ws.Range("I132").Value = 804585.9399999999
$ws.Range("K132").Value = 2413757.82
$ws.Range("M132").Value = -2411227.82
$ws.Range("H136").Value = 1258.5625
$ws.Range("I136").Value = 1152.6428
$ws.Range("K136").Value = 3457.9284
$ws.Range("M136").Value = -907.9284000000002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 20000000
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
$ws.Range("H81").Value = 1391.6666
$ws.Range("I81").Value = 883.5
$ws.Range("J81").Value = 2103.1
$ws.Range("K81").Value = 1767
$ws.Range("L81").Value = 4206.2
$ws.Range("M81").Value = -706
$ws.Range("N81").Value = -6328.2
$ws.Range("H84").Value = 1391.6666
$ws.Range("I84").Value = 883.5
$ws.Range("J84").Value = 2103.1
$ws.Range("K84").Value = 8835
$ws.Range("L84").Value = 21031
$ws.Range("M84").Value = -3531
$ws.Range("N84").Value = -31639
$ws.Range("H100").Value = 371.5
$ws.Range("I100").Value = 383.125
$ws.Range("J100").Value = 325
$ws.Range("K100").Value = 766.25
$ws.Range("L100").Value = 650
$ws.Range("M100").Value = -225.25
$ws.Range("N100").Value = -1732
$ws.Range("H126").Value = 1581.9333
$ws.Range("I126").Value = 1208.5238
$ws.Range("K126").Value = 3625.5714
$ws.Range("M126").Value = -1155.5714
$ws.Range("H132").Value = 1206.3889
$ws.Range("I132").Value = 857.25
$ws.Range("J132").Value = 3999.5
$ws.Range("K132").Value = 2571.75
$ws.Range("L132").Value = 11998.5
$ws.Range("M132").Value = -41.75
$ws.Range("N132").Value = -17058.5
